$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C29").Value = 378
$ws.Range("D29").Value = 43
$ws.Range("E29").Value = 335
$ws.Range("F29").Value = 7.401032702237521
